$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 5-7 (the FAPs->ECs, MuSCs->ECs, MuSCs->FAPs rows no longer needed
# after recomputing with the new TPM values) and shift remaining rows up.
$ws.Range("A5:T7").EntireRow.Delete()

# Row 2: ECs -> FAPs (target cluster), with recomputed TPM-based values.
$ws.Range("D2").Value = "FAPs"
$ws.Range("G2").Value = 7.292394999999999
$ws.Range("H2").Value = 21.877185
$ws.Range("I2").Value = 0.5244715940033005
$ws.Range("J2").Value = 0.5244715940033005
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.08363033333333332
$ws.Range("N2").Value = 0.250891
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.6098654246483332
$ws.Range("R2").Value = 5.488788821834999
$ws.Range("S2").Value = 0.5244715940033005
$ws.Range("T2").Value = 0.5244715940033005

# Row 3: Sending cluster ECs -> FAPs, target cluster stays FAPs, recomputed values.
$ws.Range("A3").Value = "FAPs"
$ws.Range("G3").Value = 0.327332
$ws.Range("H3").Value = 0.9819960000000001
$ws.Range("I3").Value = 0.02354183170388992
$ws.Range("J3").Value = 0.02354183170388992
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.02737488427066666
$ws.Range("R3").Value = 0.246373958436
$ws.Range("S3").Value = 0.02354183170388992
$ws.Range("T3").Value = 0.02354183170388992

# Row 4: Sending cluster FAPs -> MuSCs, target cluster ECs -> FAPs, recomputed values.
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("G4").Value = 6.284543666666667
$ws.Range("H4").Value = 18.853631
$ws.Range("I4").Value = 0.4519865742928097
$ws.Range("J4").Value = 0.4519865742928096
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.08363033333333332
$ws.Range("N4").Value = 0.250891
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.5255784816912221
$ws.Range("R4").Value = 4.730206335220999
$ws.Range("S4").Value = 0.4519865742928097
$ws.Range("T4").Value = 0.4519865742928096
